# 9th Stab - Cosmetic Changes
#
# This "MarketBeat Rank" tracker keeps one column per weekly snapshot
# (newest week right after the ticker name / current-status columns).
# A new week's worth of columns is inserted, the older week headers shift
# two slots to the right, and the new, newest column is seeded with "UN"
# (unchanged) for every analyst row except one, which picks up a
# rating-change note that gets highlighted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current newest-week header (column B) before it gets
# shifted out of the way by the column insert below.
$prevHeader = $ws.Range("B1").Value2

# Insert two new columns at C:D -- everything from the old column C onward
# (here just column C, "Jun_10") slides two columns to the right, to E.
$ws.Columns("C:D").Insert()

# The old column B header ("Jun_13") survived the insert untouched (it is
# to the left of the inserted range) but now belongs two weeks back, so it
# needs to move into the new column D; column E already holds the old
# column C header ("Jun_10") courtesy of the insert/shift itself.
$ws.Range("D1").Value = $prevHeader

# Give the two freshly inserted week columns their own headers.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Match the explicit 8.0-character width already used on the neighboring
# week columns.
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667

# Every analyst starts the new weeks as "UN" (unchanged) by default.
For ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# ValuEngine (row 11) was downgraded the week of Jun_15 -- record the
# rating change note in the new column C and highlight the cell in orange.
$c11 = $ws.Range("C11")
$c11.Value = "6/13/2018,Downgrades,Buy -> Hold,"
$c11.Interior.ColorIndex = 45
